$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "staffID"
$ws.Range("D1").Value = "role"
$ws.Range("E1").Value = "gender"
$ws.Range("F1").Value = "age"
$ws.Range("G1").Value = "branch"

# --- Data rows (2-8) ---
$ws.Range("A2").Value = "477c0c7e-9d46-4202-969d-f3dd1933a575"
$ws.Range("B2").Value = "kumar Blackmore"
$ws.Range("C2").Value = "kumarB"
$ws.Range("D2").Value = "S"
$ws.Range("E2").Value = "M"
$ws.Range("F2").Value = 32
$ws.Range("G2").Value = "NTU"

$ws.Range("A3").Value = "67136f7c-fcd0-45f1-8859-9e3d183faeb3"
$ws.Range("B3").Value = "Alexei "
$ws.Range("C3").Value = "Alexei"
$ws.Range("D3").Value = "M"
$ws.Range("E3").Value = "M"
$ws.Range("F3").Value = 25
$ws.Range("G3").Value = "NTU"

$ws.Range("A4").Value = "5cc0e578-41b6-4e7d-b6e8-5f287be3e857"
$ws.Range("B4").Value = "Tom Chan "
$ws.Range("C4").Value = "TomC"
$ws.Range("D4").Value = "M"
$ws.Range("E4").Value = "M"
$ws.Range("F4").Value = 56
$ws.Range("G4").Value = "JP"

$ws.Range("A5").Value = "6cd0a2b8-2412-4c2a-bf28-b52d043d414b"
$ws.Range("B5").Value = "Alica Ang"
$ws.Range("C5").Value = "AlicaA"
$ws.Range("D5").Value = "M"
$ws.Range("E5").Value = "F"
$ws.Range("F5").Value = 27
$ws.Range("G5").Value = "JE"

$ws.Range("A6").Value = "9234b60f-23e1-4b28-9732-c987de0605a8"
$ws.Range("B6").Value = "Mary lee"
$ws.Range("C6").Value = "MaryL"
$ws.Range("D6").Value = "S"
$ws.Range("E6").Value = "F"
$ws.Range("F6").Value = 44
$ws.Range("G6").Value = "JE"

$ws.Range("A7").Value = "22148748-4c7e-4331-8686-d2c6d3c27e22"
$ws.Range("B7").Value = "Justin Loh"
$ws.Range("C7").Value = "JustinL"
$ws.Range("D7").Value = "S"
$ws.Range("E7").Value = "M"
$ws.Range("F7").Value = 49
$ws.Range("G7").Value = "JP"

$ws.Range("A8").Value = "ec084e54-155c-4a11-b8e8-04df6cfe3c87"
$ws.Range("B8").Value = "Boss"
$ws.Range("C8").Value = "boss"
$ws.Range("D8").Value = "A"
$ws.Range("E8").Value = "F"
$ws.Range("F8").Value = 62
$ws.Range("G8").Value = "null"

# --- Hyperlinks for staff emails (order defines rId1..rId4) ---
# Adding a hyperlink with a display-text argument overwrites the cell's
# text, so immediately restore the real cell value afterwards; the
# <hyperlink display="..."/> attribute itself stays put.
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:ARVI@NTU.EDU.SG", "", "", "ARVI@NTU.EDU.SG")
$ws.Range("C6").Value = "MaryL"
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:ANWIT@NTU.EDU.SG", "", "", "ANWIT@NTU.EDU.SG")
$ws.Range("C5").Value = "AlicaA"
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:UPAM@NTU.EDU.SG", "", "", "UPAM@NTU.EDU.SG")
$ws.Range("C4").Value = "TomC"
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:OURIN@ntu.edu.sg", "", "", "OURIN@ntu.edu.sg")
$ws.Range("C3").Value = "Alexei"

# --- Carry over formatting from the columns that used to hold this data ---
# (done last so it overrides any auto-applied "Hyperlink" style from above)
# C2:C6 inherit the old "staffID-ish" look that used to live in B2:B6 (s=9)
$ws.Range("B2").Copy()
$ws.Range("C2:C6").PasteSpecial(-4122)

# C1, C7, C8 inherit the look that used to live in B1/B7/B8 (s=10)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("B7").Copy()
$ws.Range("C7:C8").PasteSpecial(-4122)

# B2:B6 now take on the s=10 look (the data that used to be there moved to C)
$ws.Range("B1").Copy()
$ws.Range("B2:B6").PasteSpecial(-4122)

# A2:A8 (the new id column) has no special styling
$ws.Range("A2:A8").Style = "Normal"

$excel.CutCopyMode = 0

# --- Selection matches the author's final view state ---
[void]$ws.Range("A2:A8").Select()
